$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.383.31'
$ws.Range('E2').Value = '  +3.19%  '
$ws.Range('D3').Value = '2.508.11'
$ws.Range('E3').Value = '  +2.47%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '110.39'
$ws.Range('E5').Value = '  +5.92%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '324.28'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.525'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +1.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.47'
$ws.Range('E10').Value = '  +10.39%  '
$ws.Range('E11').Value = '  +1.73%  '
$ws.Range('E12').Value = '  +0.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.60'
$ws.Range('E13').Value = '  +2.49%  '
$ws.Range('E14').Value = '  +2.75%  '
$ws.Range('D15').Value = '2.901.42'
$ws.Range('E15').Value = '  +2.53%  '
$ws.Range('D16').Value = '2.517.65'
$ws.Range('E16').Value = '  +2.26%  '
$ws.Range('E17').Value = '  +2.80%  '
$ws.Range('D18').Value = '47.344.08'
$ws.Range('E18').Value = '  +3.42%  '
$ws.Range('E19').Value = '  +3.71%  '
$ws.Range('E20').Value = '  +4.54%  '
$ws.Range('D21').Value = '0.0₃0944'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('E22').Value = '  +12.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.65'
$ws.Range('E23').Value = '  -0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '249.35'
$ws.Range('E24').Value = '  +1.53%  '
$ws.Range('E25').Value = '  +4.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.15'
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.29'
$ws.Range('E28').Value = '  +4.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.05'
$ws.Range('E29').Value = '  +3.78%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.66'
$ws.Range('E30').Value = '  +6.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.138'
$ws.Range('E31').Value = '  +8.45%  '
$ws.Range('E32').Value = '  +2.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.10'
$ws.Range('E33').Value = '  +1.42%  '
$ws.Range('E34').Value = '  +1.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0801'
$ws.Range('E35').Value = '  +5.73%  '
$ws.Range('E36').Value = '  +0.26%  '
$ws.Range('E37').Value = '  +6.17%  '
$ws.Range('E38').Value = '  +5.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.01'
$ws.Range('E39').Value = '  +2.99%  '
$ws.Range('E40').Value = '  +1.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '122.00'
$ws.Range('E41').Value = '  -3.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.23'
$ws.Range('E42').Value = '  -1.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.57'
$ws.Range('E43').Value = '  +3.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0300'
$ws.Range('E44').Value = '  +2.71%  '
$ws.Range('D45').Value = '2.006.58'
$ws.Range('E46').Value = '  +5.43%  '
$ws.Range('E47').Value = '  -1.79%  '
$ws.Range('E48').Value = '  -3.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.07'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.25'
$ws.Range('E50').Value = '  +6.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.42'
$ws.Range('E51').Value = '  +1.50%  '
